$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Add two new columns (G: "Codigo", H: "Construcción Iteración 1") that
#    mirror the existing nomenclature table, following the same visual
#    style as column A (header / sub-header / body banding).
# ---------------------------------------------------------------------------

# Copy column A's formatting (fills/borders) onto the new G and H columns so
# the new cells reuse the same style classes Excel already has defined.
$ws.Range("A1:A23").Copy() | Out-Null
$ws.Range("G1:G23").PasteSpecial(-4122) | Out-Null
$ws.Range("A1:A23").Copy() | Out-Null
$ws.Range("H1:H23").PasteSpecial(-4122) | Out-Null

# Populate values in the same order they were authored so the shared-string
# table grows in the same sequence as the source workbook.
$ws.Range("G1").Value = "Codigo"
$ws.Range("G2").Value = "C101"
$ws.Range("H2").Value = "Plan de Iteración"

$ws.Range("H1").Value = "Construcción Iteración 1"
$ws.Range("H3").Value = "CRUD Escenario.docx"
$ws.Range("G3").Value = "C102"
$ws.Range("G4").Value = "C103"
$ws.Range("H4").Value = "Arquitectura del Sistema"
$ws.Range("G5").Value = "C104"
$ws.Range("G6").Value = "C105"
$ws.Range("G7").Value = "C106"
$ws.Range("G8").Value = "C107"
$ws.Range("G9").Value = "C108"
$ws.Range("G10").Value = "C109"
$ws.Range("G11").Value = "C110"
$ws.Range("G12").Value = "C111"
$ws.Range("G13").Value = "C112"
$ws.Range("G14").Value = "C113"
$ws.Range("G15").Value = "C114"
$ws.Range("G16").Value = "C115"

# ---------------------------------------------------------------------------
# 2. Documentation clean-up on the existing table: two risk-report entries
#    now need to wrap their text (narrower column F below).
# ---------------------------------------------------------------------------
$ws.Range("F15:F16").WrapText = $true

# ---------------------------------------------------------------------------
# 3. Resize columns: narrow "spacer" columns for the code columns, and a
#    narrower column F now that it no longer needs to best-fit its text.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 6.163333333333343   # A  ~7
$ws.Columns.Item(2).ColumnWidth = 21.996666666666698  # B  ~22.85546875
$ws.Columns.Item(3).ColumnWidth = 6.33                # C  ~7.140625
$ws.Columns.Item(4).ColumnWidth = 25.996666666666684  # D  ~26.85546875
$ws.Columns.Item(5).ColumnWidth = 6.33                # E  ~7.140625
$ws.Columns.Item(6).ColumnWidth = 28.16333333333337   # F  29
$ws.Columns.Item(7).ColumnWidth = 6.33                # G  ~7.140625
$ws.Columns.Item(8).ColumnWidth = 32.33000000000001   # H  ~33.140625

# ---------------------------------------------------------------------------
# 4. Restore the last-used selection cursor.
# ---------------------------------------------------------------------------
$ws.Range("H9").Select() | Out-Null

Write-Host "done"
